# Append two new records (rows 35 and 36) to the "completedRegister" sheet,
# duplicating the last existing record (row 34) verbatim, then updating only
# the trailing timestamp column (BT) on each new row — mirrors the author's
# "add command line arg for check values" commit, which appended two more
# rows of the same checked case to the register.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 34 -> new row 35 (shifts old row 35+ down, none exist here).
# Copy + Insert preserves cell styles/number formats exactly (reuses the
# existing style table instead of minting new xf records).
$ws.Rows.Item(34).Copy()
$ws.Rows.Item(35).Insert()

# Duplicate row 34 again -> new row 36.
$ws.Rows.Item(34).Copy()
$ws.Rows.Item(36).Insert()

# Each appended record gets its own submission timestamp.
$ws.Range("BT35").Value = 44611.82183226852
$ws.Range("BT36").Value = 44611.82285181713
